$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: seed the shared-string table with the 8 corrected academic-year
# labels, in the exact order they should receive new shared-string indices
# (2011/12 .. 2018/19). We stash them in scratch cells far outside the
# worksheet's used range (A1:M270) -- column O (15), which (unlike columns
# D..M) has no per-column style, so writing/clearing it leaves no residue --
# then blank the scratch cells again once every real cell that needs them
# has been written. Because the engine builds the shared-string table from
# the order cells are first touched (not their row/column position), seeding
# here first locks in the desired ordering before any of the "real" column C
# updates below run.
$yearLabels = @("2011/12","2012/13","2013/14","2014/15","2015/16","2016/17","2017/18","2018/19")
$scratchCol = 15
for ($i = 0; $i -lt $yearLabels.Count; $i++) {
    $r = 1000 + $i
    $ws.Cells.Item($r, $scratchCol).Value = $yearLabels[$i]
}

# --- Step 2: fix up the malformed academic-year strings in column C.
# The original workbook stored these as "20NN/N" (e.g. "2011/2"); every row
# that has real percentage-breakdown data (as opposed to the "Fewer than 10
# enrolled" / "Not enough info" placeholder rows) gets corrected to the
# proper "20NN/NN" form (e.g. "2011/12") so the table's new search/filter
# functionality can match on it correctly.
$ws.Range("C3").Value = "2012/13"
$ws.Range("C4").Value = "2013/14"
$ws.Range("C27").Value = "2011/12"
$ws.Range("C28").Value = "2012/13"
$ws.Range("C29").Value = "2013/14"
$ws.Range("C30").Value = "2014/15"
$ws.Range("C31").Value = "2015/16"
$ws.Range("C36").Value = "2011/12"
$ws.Range("C37").Value = "2012/13"
$ws.Range("C38").Value = "2013/14"
$ws.Range("C39").Value = "2014/15"
$ws.Range("C40").Value = "2015/16"
$ws.Range("C41").Value = "2016/17"
$ws.Range("C42").Value = "2017/18"
$ws.Range("C43").Value = "2018/19"
$ws.Range("C44").Value = "2011/12"
$ws.Range("C47").Value = "2014/15"
$ws.Range("C49").Value = "2016/17"
$ws.Range("C50").Value = "2017/18"
$ws.Range("C51").Value = "2018/19"
$ws.Range("C53").Value = "2012/13"
$ws.Range("C54").Value = "2013/14"
$ws.Range("C55").Value = "2014/15"
$ws.Range("C56").Value = "2015/16"
$ws.Range("C57").Value = "2016/17"
$ws.Range("C58").Value = "2017/18"
$ws.Range("C59").Value = "2018/19"
$ws.Range("C60").Value = "2012/13"
$ws.Range("C61").Value = "2013/14"
$ws.Range("C62").Value = "2014/15"
$ws.Range("C63").Value = "2015/16"
$ws.Range("C64").Value = "2016/17"
$ws.Range("C65").Value = "2017/18"
$ws.Range("C66").Value = "2018/19"
$ws.Range("C67").Value = "2012/13"
$ws.Range("C68").Value = "2013/14"
$ws.Range("C69").Value = "2014/15"
$ws.Range("C70").Value = "2015/16"
$ws.Range("C71").Value = "2016/17"
$ws.Range("C72").Value = "2012/13"
$ws.Range("C73").Value = "2013/14"
$ws.Range("C74").Value = "2014/15"
$ws.Range("C75").Value = "2015/16"
$ws.Range("C76").Value = "2016/17"
$ws.Range("C77").Value = "2017/18"
$ws.Range("C83").Value = "2016/17"
$ws.Range("C84").Value = "2017/18"
$ws.Range("C86").Value = "2012/13"
$ws.Range("C87").Value = "2013/14"
$ws.Range("C88").Value = "2014/15"
$ws.Range("C89").Value = "2015/16"
$ws.Range("C90").Value = "2016/17"
$ws.Range("C91").Value = "2017/18"
$ws.Range("C92").Value = "2018/19"
$ws.Range("C93").Value = "2012/13"
$ws.Range("C94").Value = "2013/14"
$ws.Range("C95").Value = "2014/15"
$ws.Range("C97").Value = "2014/15"
$ws.Range("C98").Value = "2015/16"
$ws.Range("C99").Value = "2016/17"
$ws.Range("C100").Value = "2017/18"
$ws.Range("C101").Value = "2018/19"
$ws.Range("C102").Value = "2014/15"
$ws.Range("C103").Value = "2015/16"
$ws.Range("C104").Value = "2016/17"
$ws.Range("C105").Value = "2017/18"
$ws.Range("C106").Value = "2017/18"
$ws.Range("C107").Value = "2018/19"
$ws.Range("C108").Value = "2017/18"
$ws.Range("C110").Value = "2016/17"
$ws.Range("C111").Value = "2017/18"
$ws.Range("C112").Value = "2018/19"
$ws.Range("C113").Value = "2017/18"
$ws.Range("C114").Value = "2018/19"
$ws.Range("C115").Value = "2017/18"
$ws.Range("C116").Value = "2018/19"
$ws.Range("C117").Value = "2017/18"
$ws.Range("C118").Value = "2018/19"
$ws.Range("C119").Value = "2018/19"
$ws.Range("C122").Value = "2014/15"
$ws.Range("C123").Value = "2011/12"
$ws.Range("C125").Value = "2013/14"
$ws.Range("C127").Value = "2015/16"
$ws.Range("C130").Value = "2012/13"
$ws.Range("C131").Value = "2013/14"
$ws.Range("C132").Value = "2014/15"
$ws.Range("C133").Value = "2015/16"
$ws.Range("C134").Value = "2016/17"
$ws.Range("C135").Value = "2017/18"
$ws.Range("C136").Value = "2018/19"
$ws.Range("C138").Value = "2012/13"
$ws.Range("C139").Value = "2013/14"
$ws.Range("C140").Value = "2014/15"
$ws.Range("C141").Value = "2015/16"
$ws.Range("C142").Value = "2016/17"
$ws.Range("C143").Value = "2017/18"
$ws.Range("C144").Value = "2018/19"
$ws.Range("C148").Value = "2011/12"
$ws.Range("C152").Value = "2011/12"
$ws.Range("C153").Value = "2012/13"
$ws.Range("C156").Value = "2015/16"
$ws.Range("C157").Value = "2016/17"
$ws.Range("C158").Value = "2017/18"
$ws.Range("C159").Value = "2018/19"
$ws.Range("C167").Value = "2011/12"
$ws.Range("C168").Value = "2012/13"
$ws.Range("C170").Value = "2014/15"
$ws.Range("C172").Value = "2016/17"
$ws.Range("C174").Value = "2018/19"
$ws.Range("C178").Value = "2018/19"
$ws.Range("C179").Value = "2011/12"
$ws.Range("C181").Value = "2013/14"
$ws.Range("C183").Value = "2015/16"
$ws.Range("C193").Value = "2011/12"
$ws.Range("C194").Value = "2011/12"
$ws.Range("C195").Value = "2012/13"
$ws.Range("C196").Value = "2013/14"
$ws.Range("C197").Value = "2014/15"
$ws.Range("C198").Value = "2015/16"
$ws.Range("C199").Value = "2016/17"
$ws.Range("C200").Value = "2017/18"
$ws.Range("C201").Value = "2018/19"
$ws.Range("C202").Value = "2011/12"
$ws.Range("C203").Value = "2012/13"
$ws.Range("C204").Value = "2013/14"
$ws.Range("C205").Value = "2014/15"
$ws.Range("C206").Value = "2015/16"
$ws.Range("C207").Value = "2016/17"
$ws.Range("C208").Value = "2017/18"
$ws.Range("C209").Value = "2018/19"
$ws.Range("C210").Value = "2011/12"
$ws.Range("C211").Value = "2012/13"
$ws.Range("C212").Value = "2013/14"
$ws.Range("C213").Value = "2014/15"
$ws.Range("C214").Value = "2015/16"
$ws.Range("C215").Value = "2016/17"
$ws.Range("C216").Value = "2017/18"
$ws.Range("C217").Value = "2018/19"
$ws.Range("C218").Value = "2011/12"
$ws.Range("C219").Value = "2012/13"
$ws.Range("C220").Value = "2013/14"
$ws.Range("C221").Value = "2014/15"
$ws.Range("C222").Value = "2015/16"
$ws.Range("C223").Value = "2016/17"
$ws.Range("C224").Value = "2017/18"
$ws.Range("C225").Value = "2018/19"
$ws.Range("C226").Value = "2011/12"
$ws.Range("C227").Value = "2012/13"
$ws.Range("C228").Value = "2013/14"
$ws.Range("C229").Value = "2014/15"
$ws.Range("C230").Value = "2015/16"
$ws.Range("C231").Value = "2016/17"
$ws.Range("C232").Value = "2017/18"
$ws.Range("C233").Value = "2018/19"
$ws.Range("C234").Value = "2011/12"
$ws.Range("C235").Value = "2012/13"
$ws.Range("C236").Value = "2013/14"
$ws.Range("C237").Value = "2014/15"
$ws.Range("C238").Value = "2015/16"
$ws.Range("C239").Value = "2016/17"
$ws.Range("C240").Value = "2017/18"
$ws.Range("C241").Value = "2018/19"
$ws.Range("C242").Value = "2011/12"
$ws.Range("C243").Value = "2012/13"
$ws.Range("C244").Value = "2013/14"
$ws.Range("C245").Value = "2014/15"
$ws.Range("C246").Value = "2015/16"
$ws.Range("C247").Value = "2016/17"
$ws.Range("C248").Value = "2017/18"
$ws.Range("C249").Value = "2018/19"
$ws.Range("C250").Value = "2011/12"
$ws.Range("C251").Value = "2012/13"
$ws.Range("C252").Value = "2013/14"
$ws.Range("C253").Value = "2014/15"
$ws.Range("C254").Value = "2015/16"
$ws.Range("C255").Value = "2016/17"
$ws.Range("C256").Value = "2017/18"
$ws.Range("C259").Value = "2012/13"
$ws.Range("C261").Value = "2014/15"
$ws.Range("C262").Value = "2015/16"
$ws.Range("C263").Value = "2016/17"
$ws.Range("C264").Value = "2017/18"
$ws.Range("C265").Value = "2018/19"
$ws.Range("C266").Value = "2011/12"
$ws.Range("C267").Value = "2015/16"
$ws.Range("C268").Value = "2016/17"
$ws.Range("C269").Value = "2017/18"
$ws.Range("C270").Value = "2018/19"

# --- Step 3: clear the scratch cells now that the real cells above hold
# their own references to these strings.
for ($i = 0; $i -lt $yearLabels.Count; $i++) {
    $r = 1000 + $i
    $ws.Cells.Item($r, $scratchCol).Value = ""
}

# --- Step 4: move the active selection to C49 (matches the saved cursor
# position recorded in the workbook after this edit).
$ws.Range("C49").Select()
